# Traditional Chinese -> Simplified Chinese translation update for
# "Email 10-1 [TEMPLATE] Partner email - thank you email (without photos).docx"
#
# Strategy: use Find/Replace (wdReplaceAll) against $d.Content for every
# run whose *entire* text node is being swapped for new text - this keeps
# each edit inside the boundaries of a single <w:r>, so the run's own
# formatting (rPr) is left untouched. For the one spot where only a
# single character in the middle of a run needs to change (the full-width
# comma between [CITY] and [COUNTRY]) we locate it precisely via
# Find + Collapse + MoveEnd so the neighbouring highlighted runs for
# [CITY]/[COUNTRY] are not touched.

$d = $word.ActiveDocument

function Replace-All($old, $new) {
    $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# Language list at the top of the doc (hyperlink run + following run).
Replace-All "英文" "英语"
Replace-All " / 葡萄牙文 / 法文 / 泰文 / 越南文 / 西班牙文" " / 葡萄牙语 / 法语 / 泰语 / 越南语 / 西班牙语"

# Summary table: "簡介" (Description) / its body / "目標受眾" (Target
# audience) / its body.
Replace-All "簡介" "简介"
Replace-All "寄給參加活動的合作夥伴的電子郵件。 此電子郵件將包括照片畫廊，將通過 customer.io 發送。" `
            "一封发送给参加活动的合作伙伴的邮件。 这封邮件将包含一个照片画廊，将通过 customer.io 发送。"
Replace-All "目標受眾" "目标受众"
Replace-All "活動參加者" "活动参与者"

# "Subject: " line.
Replace-All "Subject: " "主题: "
Replace-All "感謝您參加 " "感谢您参加 "

# Heading "You made our event a success! ..."
Replace-All "您使我們的活動圓滿成功！ 🎉" "您使我们的活动圆满成功！ 🎉"

# Greeting line: "[PARTNER NAME]" placeholder + trailing punctuation run.
Replace-All "[PARTNER NAME]" "[合作伙伴姓名]"
Replace-All "， " ", "

# "Thank you for attending [EVENT NAME] in [CITY], [COUNTRY]." paragraph.
Replace-All " 於 " " 于 "
Replace-All "。 希望您度過了愉快的時光，很高興認識您！" "。 希望您玩得开心，很高兴认识您！"

# The single full-width comma that sits in its own run between [CITY]
# and [COUNTRY] - touch only that run, not the highlighted placeholders
# on either side of it.
$cityRng = $d.Content
$cityRng.Find.Execute("[CITY]") | Out-Null
$cityRng.Collapse(0)
$cityRng.MoveEnd(1, 1)
$cityRng.Text = ", "

# "To browse the conference/workshop/..." paragraph (split across three
# runs around the comment anchors).
Replace-All "如需瀏覽會議/" "如需浏览 "
Replace-All "研討會/聯盟之旅" "会议/研讨会/联盟之旅"
Replace-All "的照片和精彩片段，並隨時了解我們為您舉辦的最新活動和計劃，請關注我們的社交媒體帳戶：" `
            " 的照片和精彩片段，并随时了解我们举办的最新活动和计划，请关注我们："

# Closing line.
Replace-All "希望這次活動能給您們帶來和我們一樣的啟發，讓我們繼續共同成長！" "希望这次活动能给您们带来和我们一样的启发，让我们继续共同成长！"

# Comment text ("Choose one of these").  Comments in this runtime are not
# reachable through Range.Find (the Comment.Range anchor is not scoped to
# the comments story), so fall back to assigning Range.Text directly,
# which is the one operation that reliably reaches the comment body.
$d.Comments.Item(1).Range.Text = "选择其中之一"
